$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 411 (existing rows 411..442 shift down to 413..444)
$ws.Rows.Item(411).Resize(2).Insert() | Out-Null

# New row 411: Lane Late / Primera
$ws.Cells.Item(411, 1).Value = 3
$ws.Cells.Item(411, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(411, 3).Value = "Coquimbo"
$ws.Cells.Item(411, 4).Value = 44461
$ws.Cells.Item(411, 5).Value = 5
$ws.Cells.Item(411, 6).Value = "Fruta"
$ws.Cells.Item(411, 7).Value = 100102
$ws.Cells.Item(411, 8).Value = "Cítricos"
$ws.Cells.Item(411, 9).Value = 100102005
$ws.Cells.Item(411, 10).Value = "Naranja"
$ws.Cells.Item(411, 11).Value = "Lane Late"
$ws.Cells.Item(411, 12).Value = "Primera"
$ws.Cells.Item(411, 13).Value = 148
$ws.Cells.Item(411, 14).Value = 4500
$ws.Cells.Item(411, 15).Value = 5000
$ws.Cells.Item(411, 16).Value = 4736
$ws.Cells.Item(411, 17).Value = "$/malla 13 kilos"
$ws.Cells.Item(411, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(411, 19).Value = 364
$ws.Cells.Item(411, 20).Value = 13

# New row 412: Lane Late / Segunda
$ws.Cells.Item(412, 1).Value = 3
$ws.Cells.Item(412, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(412, 3).Value = "Coquimbo"
$ws.Cells.Item(412, 4).Value = 44461
$ws.Cells.Item(412, 5).Value = 5
$ws.Cells.Item(412, 6).Value = "Fruta"
$ws.Cells.Item(412, 7).Value = 100102
$ws.Cells.Item(412, 8).Value = "Cítricos"
$ws.Cells.Item(412, 9).Value = 100102005
$ws.Cells.Item(412, 10).Value = "Naranja"
$ws.Cells.Item(412, 11).Value = "Lane Late"
$ws.Cells.Item(412, 12).Value = "Segunda"
$ws.Cells.Item(412, 13).Value = 163
$ws.Cells.Item(412, 14).Value = 3500
$ws.Cells.Item(412, 15).Value = 4000
$ws.Cells.Item(412, 16).Value = 3761
$ws.Cells.Item(412, 17).Value = "$/malla 13 kilos"
$ws.Cells.Item(412, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(412, 19).Value = 289
$ws.Cells.Item(412, 20).Value = 13
